# This script applies a cyclic shift to the species-observation rows 64-68
# on the active worksheet: the record that used to live on row 68 moves up
# to row 64, and the records that used to be on rows 64-67 each shift down
# by one row (to rows 65-68 respectively). Only the columns that actually
# vary between these rows are touched: A, B, D, E, F, G, H, Q, R, Z, AB.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current ("before") values for the varying columns of rows 64-68.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")
$rows = 64..68

$before = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowData
}

# Compute the "after" mapping: row 64 gets what row 68 had, and rows
# 65-68 get what rows 64-67 had (a rotation of the block).
$after = @{}
$after[64] = $before[68]
$after[65] = $before[64]
$after[66] = $before[65]
$after[67] = $before[66]
$after[68] = $before[67]

# Write the new values back.
foreach ($r in $rows) {
    $rowData = $after[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rowData[$c]
    }
}
